# FW_UI_0000.xlsx — "8.16 Sprint 2 cases"
#
# Adds new Feeds-related test cases (FW_UI_0014..FW_UI_0021) to the
# RunTest sheet, re-striping the alternating row-banding styles (style
# index 2 / style index 5) for rows 6-23, and drops the now-unused last
# blank row (24), which moves dimension A1:F24 -> A1:F23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style changes -----------------------------------------------------
# Row 2 is entirely "style 2" and row 3 is entirely "style 5" in both the
# original and final workbook, so they're safe, stable format sources for
# every PasteSpecial(Formats) below regardless of execution order.
$style2Src = $ws.Range("A2")
$style5Src = $ws.Range("A3")

# Row 6 (FW_UI_0004): B flips 2 -> 5
$style5Src.Copy()
$ws.Range("B6").PasteSpecial(-4122)

# Row 7 (FW_UI_0005): B flips 5 -> 2
$style2Src.Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Row 8 (FW_UI_0006): B flips 2 -> 5
$style5Src.Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Row 12 (FW_UI_0010): B flips 2 -> 5
$style5Src.Copy()
$ws.Range("B12").PasteSpecial(-4122)

# Row 13 (FW_UI_0011): B flips 5 -> 2
$style2Src.Copy()
$ws.Range("B13").PasteSpecial(-4122)

# Row 14 (FW_UI_0012): A,C,D,E,F flip 5 -> 2 (B stays 5)
$style2Src.Copy()
$ws.Range("A14").PasteSpecial(-4122)
$style2Src.Copy()
$ws.Range("C14:F14").PasteSpecial(-4122)

# Row 15 (FW_UI_0013, existing case, re-striped + gains a Parameters value)
$style5Src.Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

# Row 16 (new: FW_UI_0014 / VerifyFeedinHeadline)
$style2Src.Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)

# Row 17 (new: FW_UI_0015 / VerifyFeedRemoval)
$style5Src.Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)

# Row 18 (new: FW_UI_0016 / VerifyFeedDeselection)
$style2Src.Copy()
$ws.Range("A18").PasteSpecial(-4122)
$style2Src.Copy()
$ws.Range("C18:F18").PasteSpecial(-4122)

# Row 19 (new: FW_UI_0017 / VerifyFeedReselection)
$style5Src.Copy()
$ws.Range("A19").PasteSpecial(-4122)
$style5Src.Copy()
$ws.Range("D19:F19").PasteSpecial(-4122)

# Row 20 (new: FW_UI_0018 / VerifyFeedRelaunch)
$style2Src.Copy()
$ws.Range("A20").PasteSpecial(-4122)
$style2Src.Copy()
$ws.Range("C20:F20").PasteSpecial(-4122)

# Row 21 (new: FW_UI_0019 / VerifyFeedinHeadline, torn-out-tab variant)
$style5Src.Copy()
$ws.Range("A21:F21").PasteSpecial(-4122)

# Row 22 (new placeholder: FW_UI_0020)
$style2Src.Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)

# Row 23 (new placeholder: FW_UI_0021)
$style5Src.Copy()
$ws.Range("A23:F23").PasteSpecial(-4122)

# --- Value changes -------------------------------------------------------
$ws.Range('A15').Value2 = 'No'
$ws.Range('D15').Value2 = 'Australia,sydney.newsroom,Denmark,copenhagen.newsroom'

$ws.Range('A16').Value2 = 'No'
$ws.Range('B16').Value2 = 'FW_UI_0014'
$ws.Range('C16').Value2 = 'VerifyFeedinHeadline'
$ws.Range('D16').Value2 = 'Australia,sydney.newsroom,Denmark,copenhagen.newsroom,normal'
$ws.Range('E16').Value2 = 'Feeds233961'
$ws.Range('F16').Value2 = 'To verify that user is able to view the releases in FW UI based on selected feeds'

$ws.Range('A17').Value2 = 'No'
$ws.Range('B17').Value2 = 'FW_UI_0015'
$ws.Range('C17').Value2 = 'VerifyFeedRemoval'
$ws.Range('D17').Value2 = 'Australia,sydney.newsroom'
$ws.Range('E17').Value2 = 'Feeds233961'
$ws.Range('F17').Value2 = 'To verify that  user can remove selected feeds Webui preferences'

$ws.Range('A18').Value2 = 'No'
$ws.Range('B18').Value2 = 'FW_UI_0016'
$ws.Range('C18').Value2 = 'VerifyFeedDeselection'
$ws.Range('D18').Value2 = 'Australia,sydney.newsroom'
$ws.Range('E18').Value2 = 'Feeds233961'
$ws.Range('F18').Value2 = 'To verify that user is able to  deselect feeds  from Feeds dropdown'

$ws.Range('A19').Value2 = 'No'
$ws.Range('B19').Value2 = 'FW_UI_0017'
$ws.Range('C19').Value2 = 'VerifyFeedReselection'
$ws.Range('D19').Value2 = 'Australia,sydney.newsroom'
$ws.Range('E19').Value2 = 'Feeds233961'
$ws.Range('F19').Value2 = 'To verify that user is able to select feeds  from Feeds dropdown'

$ws.Range('A20').Value2 = 'No'
$ws.Range('B20').Value2 = 'FW_UI_0018'
$ws.Range('C20').Value2 = 'VerifyFeedRelaunch'
$ws.Range('D20').Value2 = 'Australia,sydney.newsroom'
$ws.Range('E20').Value2 = 'Feeds233961'
$ws.Range('F20').Value2 = 'To verify that user is able to view the releases in FW UI based on selected feeds in relaunched LE'

$ws.Range('A21').Value2 = 'Yes'
$ws.Range('B21').Value2 = 'FW_UI_0019'
$ws.Range('C21').Value2 = 'VerifyFeedinHeadline'
$ws.Range('D21').Value2 = 'Australia,sydney.newsroom,Denmark,copenhagen.newsroom,torn'
$ws.Range('E21').Value2 = 'Feeds233961'
$ws.Range('F21').Value2 = 'To verify that user is able to view the releases in FW UI based on selected feeds in torn out tab'

$ws.Range('B22').Value2 = 'FW_UI_0020'
$ws.Range('B23').Value2 = 'FW_UI_0021'

# --- Drop the now-unused trailing blank row (24), dimension becomes A1:F23
$ws.Rows("24:24").Delete()

# --- Tidy up: clear marching-ants clipboard marker, set final selection --
$excel.CutCopyMode = 0
[void]$ws.Range("A19").Select()
